# Top10Mapping-es.xlsx: update the 2017/2021 mapping table text, drop two
# stray leftover connector shapes, retitle a couple of columns, and move
# the active selection — mirrors the gh-pages deploy commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Cell text updates (column C: 2017 list translated to English;
#    column D "(Nuevo)" -> "(Nueva)"; column E: several 2021 wording
#    tweaks).
# ---------------------------------------------------------------------
$ws.Range("C7").Value  = "A01:2017-Injection"
$ws.Range("C8").Value  = "A02:2017-Broken Authentication"
$ws.Range("C9").Value  = "A03:2017-Sensitive Data Exposure"
$ws.Range("C10").Value = "A04:2017-XML External Entities (XXE)"
$ws.Range("C11").Value = "A05:2017-Broken Access Control"
$ws.Range("C12").Value = "A06:2017-Security Misconfiguration"
$ws.Range("C13").Value = "A07:2017-Cross-Site Scripting (XSS)"
$ws.Range("C14").Value = "A08:2017-Insecure Deserialization"
$ws.Range("C15").Value = "A09:2017-Using Components with Known Vulnerabilities"
$ws.Range("C16").Value = "A10:2017-Insufficient Logging & Monitoring"

$ws.Range("D10").Value = "(Nueva)"
$ws.Range("D14").Value = "(Nueva)"
$ws.Range("D16").Value = "(Nueva)"

$ws.Range("E8").Value  = "A02:2021-Fallas Criptográficas"
$ws.Range("E11").Value = "A05:2021-Configuración de Seguridad Incorrecta"
$ws.Range("E12").Value = "A06:2021-Componentes Vulnerables y Desactualizados"
$ws.Range("E15").Value = "A09:2021-Fallas en el Registro y Monitoreo*"
$ws.Range("E16").Value = "A10:2021-Falsificación de Solicitudes del Lado del Servidor (SSRF)*"

# ---------------------------------------------------------------------
# 2. Remove the two duplicate/stray "Straight Arrow Connector" shapes
#    (11 and 12) left over in the drawing layer.
# ---------------------------------------------------------------------
$ws.Shapes.Item("Straight Arrow Connector 11").Delete()
$ws.Shapes.Item("Straight Arrow Connector 12").Delete()

# ---------------------------------------------------------------------
# 3. Column C width change (54.33 -> 48.5 characters).
# ---------------------------------------------------------------------
$ws.Columns.Item(3).ColumnWidth = 47.666666666666664

# ---------------------------------------------------------------------
# 4. Move the active selection to E23.
# ---------------------------------------------------------------------
$ws.Range("E23").Select() | Out-Null

# ---------------------------------------------------------------------
# 5. Reposition/resize the workbook window (matches the saved
#    bookViews/workbookView geometry on the author's machine).
# ---------------------------------------------------------------------
$win = $wb.Windows.Item(1)
$win.Left = -6100
$win.Top = -21100
$win.Width = 25600
$win.Height = 21100
